# Updated cryptos list (prices + 1h volume %) for the "cryptos" sheet.
# Numeric-looking price strings are apostrophe-prefixed so Excel stores them
# as text (matching the workbook's original inline-string / text layout)
# instead of silently re-typing them as numbers (which would drop trailing
# zeros / switch to scientific notation for very small values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.741.73"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.658.27"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'518.75"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'146.39"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.577"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "2.671.01"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'6.34"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "3.124.84"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "59.759.88"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.831.54"
$ws.Range("E16").Value = "  +6.21%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'21.23"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000138"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "'350.58"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'10.39"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "'6.29"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "'62.99"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "'0.417"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "0.0₃0813"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").Value = "'7.17"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'6.57"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "'1.59"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'18.96"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'150.44"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").Value = "'4.10"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("E36").Value = "  -11.38%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "'0.875"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "'1.51"
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("D40").Value = "'36.85"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'3.70"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").Value = "'282.32"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "'19.81"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'0.606"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "2.088.15"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("D48").Value = "'0.0534"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "'0.0233"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "'4.72"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +0.36%  "
